$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 293, shifting the existing rows 293-308 down to 297-312.
$ws.Rows("293:296").Insert()

# Populate the 4 newly inserted rows (293-296) with the new weekly price block
# (date 45041), following the same layout as the other weekly blocks.

# Row 293 - Especial
$ws.Range("A293").Value = 1
$ws.Range("B293").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C293").Value = "Arica y Parinacota"
$ws.Range("D293").Value = 45041
$ws.Range("E293").Value = 15
$ws.Range("F293").Value = "Fruta"
$ws.Range("G293").Value = 100108
$ws.Range("H293").Value = "Tropicales y subtropicales"
$ws.Range("I293").Value = 100108005
$ws.Range("J293").Value = "Piña"
$ws.Range("K293").Value = "Caramelo"
$ws.Range("L293").Value = "Especial"
$ws.Range("M293").Value = 250
$ws.Range("N293").Value = 19000
$ws.Range("O293").Value = 20000
$ws.Range("P293").Value = 19600
$ws.Range("Q293").Value = "$/caja 10 unidades"
$ws.Range("R293").Value = "Ecuador"
$ws.Range("S293").Value = 1960
$ws.Range("T293").Value = 10

# Row 294 - Primera
$ws.Range("A294").Value = 1
$ws.Range("B294").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C294").Value = "Arica y Parinacota"
$ws.Range("D294").Value = 45041
$ws.Range("E294").Value = 15
$ws.Range("F294").Value = "Fruta"
$ws.Range("G294").Value = 100108
$ws.Range("H294").Value = "Tropicales y subtropicales"
$ws.Range("I294").Value = 100108005
$ws.Range("J294").Value = "Piña"
$ws.Range("K294").Value = "Caramelo"
$ws.Range("L294").Value = "Primera"
$ws.Range("M294").Value = 220
$ws.Range("N294").Value = 19000
$ws.Range("O294").Value = 20000
$ws.Range("P294").Value = 19545
$ws.Range("Q294").Value = "$/caja 12 unidades"
$ws.Range("R294").Value = "Ecuador"
$ws.Range("S294").Value = 1629
$ws.Range("T294").Value = 12

# Row 295 - Segunda
$ws.Range("A295").Value = 1
$ws.Range("B295").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C295").Value = "Arica y Parinacota"
$ws.Range("D295").Value = 45041
$ws.Range("E295").Value = 15
$ws.Range("F295").Value = "Fruta"
$ws.Range("G295").Value = 100108
$ws.Range("H295").Value = "Tropicales y subtropicales"
$ws.Range("I295").Value = 100108005
$ws.Range("J295").Value = "Piña"
$ws.Range("K295").Value = "Caramelo"
$ws.Range("L295").Value = "Segunda"
$ws.Range("M295").Value = 200
$ws.Range("N295").Value = 19000
$ws.Range("O295").Value = 20000
$ws.Range("P295").Value = 19500
$ws.Range("Q295").Value = "$/caja 14 unidades"
$ws.Range("R295").Value = "Ecuador"
$ws.Range("S295").Value = 1393
$ws.Range("T295").Value = 14

# Row 296 - Tercera
$ws.Range("A296").Value = 1
$ws.Range("B296").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C296").Value = "Arica y Parinacota"
$ws.Range("D296").Value = 45041
$ws.Range("E296").Value = 15
$ws.Range("F296").Value = "Fruta"
$ws.Range("G296").Value = 100108
$ws.Range("H296").Value = "Tropicales y subtropicales"
$ws.Range("I296").Value = 100108005
$ws.Range("J296").Value = "Piña"
$ws.Range("K296").Value = "Caramelo"
$ws.Range("L296").Value = "Tercera"
$ws.Range("M296").Value = 200
$ws.Range("N296").Value = 19000
$ws.Range("O296").Value = 20000
$ws.Range("P296").Value = 19500
$ws.Range("Q296").Value = "$/caja 16 unidades"
$ws.Range("R296").Value = "Ecuador"
$ws.Range("S296").Value = 1219
$ws.Range("T296").Value = 16
